$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "30.433.05"
$r.Style = "Normal"
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "1.851.82"
$r.Style = "Normal"
$ws.Range("E3").Value = "  +1.23%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "233.40"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("E6").Value = "  +0.01%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.4752"
$r.Style = "Normal"
$ws.Range("E7").Value = "  +2.76%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.2750"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +2.21%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.06327"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +1.53%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "17.60"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +9.64%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "1.865.26"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +1.89%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.07460"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +1.35%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "4.948"
$r.Style = "Normal"
$ws.Range("E13").Value = "  +1.23%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "84.81"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +2.29%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.6252"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +1.27%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "30.397.01"
$r.Style = "Normal"
$ws.Range("E16").Value = "  +1.21%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "246.94"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +8.75%  "
$ws.Range("E18").Value = "  +0.01%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "12.67"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +2.78%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "0.000007325"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("E21").Value = "  +0.02%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "4.905"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +1.64%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "5.907"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +1.78%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "164.77"
$r.Style = "Normal"
$ws.Range("E24").Value = "  -0.17%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "9.098"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -0.01%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "17.98"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +1.75%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "1.872"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +1.88%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "0.1030"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +1.97%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "1.348"
$r.Style = "Normal"
$ws.Range("E29").Value = "  -1.39%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "4.035"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +0.26%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "3.821"
$r.Style = "Normal"
$ws.Range("E31").Value = "  +2.20%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "0.04843"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +1.53%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "1.131"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +0.97%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "0.6973"
$r.Style = "Normal"
$ws.Range("E34").Value = "  +0.13%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "2.710"
$r.Style = "Normal"
$ws.Range("E35").Value = "  +0.87%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.01899"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +5.36%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "2.681"
$r.Style = "Normal"
$ws.Range("E37").Value = "  +3.00%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.8783"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -1.28%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "1.994"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +4.40%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "106.82"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +4.46%  "
$ws.Range("E41").Value = "  +0.06%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.4054"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +1.99%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "5.506"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +0.69%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "7.179"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +4.11%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "63.23"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +6.92%  "
$ws.Range("E46").Value = "  +1.09%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "33.72"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +4.05%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "8.527"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +1.87%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.05502"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -0.31%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "1.350"
$r.Style = "Normal"
$ws.Range("E50").Value = "  +0.10%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.3681"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +1.81%  "
